# This script applies a 3-row cyclic rotation of data in rows 4-6 of the
# "Artfynd" sheet:
#   new row 4 <- old row 5
#   new row 5 <- old row 6
#   new row 6 <- old row 4
#
# Rather than moving cell-by-cell, we snapshot the values of the fields
# that actually differ between the rows (the rest of the row content -
# P, Q, R, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY - is
# identical across rows 4-6 and is left untouched), then write the
# rotated values back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the ".Value" getter is unreliable in this runtime (it returns a
# reflection description string instead of the actual cell value), so we
# read via ".Value2" instead, which works correctly. Writing via ".Value"
# works fine.

# --- snapshot old values -------------------------------------------------

$oldA4 = $ws.Range("A4").Value2
$oldB4 = $ws.Range("B4").Value2
$oldD4 = $ws.Range("D4").Value2
$oldE4 = $ws.Range("E4").Value2
$oldF4 = $ws.Range("F4").Value2
$oldG4 = $ws.Range("G4").Value2
$oldH4 = $ws.Range("H4").Value2
$oldS4 = $ws.Range("S4").Value2
$oldAC4 = $ws.Range("AC4").Value2
$oldAI4 = $ws.Range("AI4").Value2

$oldA5 = $ws.Range("A5").Value2
$oldB5 = $ws.Range("B5").Value2
$oldD5 = $ws.Range("D5").Value2
$oldE5 = $ws.Range("E5").Value2
$oldF5 = $ws.Range("F5").Value2
$oldG5 = $ws.Range("G5").Value2
$oldH5 = $ws.Range("H5").Value2
$oldS5 = $ws.Range("S5").Value2
$oldAC5 = $ws.Range("AC5").Value2
$oldAI5 = $ws.Range("AI5").Value2

$oldA6 = $ws.Range("A6").Value2
$oldB6 = $ws.Range("B6").Value2
$oldD6 = $ws.Range("D6").Value2
$oldE6 = $ws.Range("E6").Value2
$oldF6 = $ws.Range("F6").Value2
$oldG6 = $ws.Range("G6").Value2
$oldH6 = $ws.Range("H6").Value2
$oldS6 = $ws.Range("S6").Value2
$oldAC6 = $ws.Range("AC6").Value2
$oldAI6 = $ws.Range("AI6").Value2

# --- write new row 4 (= old row 5) ---------------------------------------

$ws.Range("A4").Value = $oldA5
$ws.Range("B4").Value = $oldB5
$ws.Range("D4").Value = $oldD5
$ws.Range("E4").Value = $oldE5
$ws.Range("F4").Value = $oldF5
$ws.Range("G4").Value = $oldG5
$ws.Range("H4").Value = $oldH5
$ws.Range("S4").Value = $oldS5
$ws.Range("AC4").ClearContents()
$ws.Range("AI4").ClearContents()

# --- write new row 5 (= old row 6) ---------------------------------------

$ws.Range("A5").Value = $oldA6
$ws.Range("B5").Value = $oldB6
$ws.Range("D5").Value = $oldD6
$ws.Range("E5").Value = $oldE6
$ws.Range("F5").Value = $oldF6
$ws.Range("G5").Value = $oldG6
$ws.Range("H5").Value = $oldH6
$ws.Range("S5").Value = $oldS6
$ws.Range("AC5").ClearContents()
$ws.Range("AI5").ClearContents()

# --- write new row 6 (= old row 4) ---------------------------------------

$ws.Range("A6").Value = $oldA4
$ws.Range("B6").Value = $oldB4
$ws.Range("D6").Value = $oldD4
$ws.Range("E6").Value = $oldE4
$ws.Range("F6").Value = $oldF4
$ws.Range("G6").Value = $oldG4
$ws.Range("H6").Value = $oldH4
$ws.Range("S6").Value = $oldS4
$ws.Range("AC6").Value = $oldAC4
$ws.Range("AI6").Value = $oldAI4
